$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matching the source workbook, which stores
# every Price/Volume cell as inline text) by pre-setting a Text number format.
$textCells = @("D4", "D5", "D7", "D9", "D11", "D12", "D13", "D14", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Updated Price (D) / Volume(1h) (E) values scraped by the cron job.
$updates = @{
    "D2" = "30.296.97"
    "E2" = "  +0.08%  "
    "D3" = "1.867.93"
    "E3" = "  -0.66%  "
    "D4" = "1.001"
    "E4" = "  +0.17%  "
    "D5" = "237.40"
    "E5" = "  +0.43%  "
    "E6" = "  +0.22%  "
    "D7" = "0.4791"
    "E7" = "  -1.14%  "
    "E8" = "  -2.82%  "
    "D9" = "0.06496"
    "E9" = "  -1.53%  "
    "D10" = "1.858.52"
    "E10" = "  -1.15%  "
    "D11" = "0.07434"
    "E11" = "  +1.45%  "
    "D12" = "16.45"
    "E12" = "  -2.60%  "
    "D13" = "5.059"
    "E13" = "  -1.64%  "
    "D14" = "87.74"
    "E14" = "  +0.42%  "
    "E15" = "  -0.70%  "
    "D16" = "30.276.24"
    "E16" = "  +0.10%  "
    "E17" = "  -0.86%  "
    "E18" = "  +0.08%  "
    "D19" = "0.000007568"
    "E19" = "  -2.38%  "
    "D20" = "2.109.77"
    "E20" = "  -0.76%  "
    "D21" = "1.001"
    "E21" = "  +0.19%  "
    "D22" = "5.268"
    "E22" = "  -2.82%  "
    "D23" = "218.53"
    "E23" = "  +11.80%  "
    "D24" = "6.144"
    "E24" = "  -0.05%  "
    "D25" = "9.295"
    "E25" = "  +0.28%  "
    "D26" = "167.58"
    "E26" = "  +2.36%  "
    "D27" = "18.36"
    "E27" = "  +1.32%  "
    "D28" = "1.967"
    "E28" = "  +2.38%  "
    "E29" = "  +1.34%  "
    "D30" = "0.09358"
    "E30" = "  +2.19%  "
    "D31" = "4.289"
    "E31" = "  +0.20%  "
    "D32" = "4.009"
    "E32" = "  -0.36%  "
    "D33" = "0.05039"
    "E33" = "  -0.72%  "
    "D34" = "1.201"
    "E34" = "  +7.71%  "
    "D35" = "0.7458"
    "E35" = "  +3.99%  "
    "D36" = "2.710"
    "D37" = "0.01815"
    "E37" = "  +1.91%  "
    "D38" = "2.610"
    "E38" = "  -1.07%  "
    "D39" = "2.070"
    "E39" = "  +1.20%  "
    "D40" = "0.9047"
    "E40" = "  -1.50%  "
    "E41" = "  +2.22%  "
    "D42" = "106.51"
    "E42" = "  +0.70%  "
    "E43" = "  +0.31%  "
    "D44" = "0.4247"
    "E44" = "  -0.94%  "
    "D45" = "7.352"
    "E45" = "  -1.06%  "
    "D46" = "0.1279"
    "E46" = "  -2.22%  "
    "D47" = "63.57"
    "E47" = "  -2.18%  "
    "D48" = "1.468"
    "E48" = "  -3.26%  "
    "D49" = "8.894"
    "E49" = "  +0.48%  "
    "D50" = "33.52"
    "E50" = "  -0.91%  "
    "D51" = "0.05622"
    "E51" = "  -2.11%  "
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output "Updated $($updates.Count) cells"
